# Reprogramacion de credito empresarial - actualizar datos de integracion
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cliente (A2) y Pagare (G2) ya estan formateados como texto (numFmt "Text"),
# por lo que basta con asignar el nuevo valor como cadena.
$ws.Range("A2").Value = "24681769"
$ws.Range("G2").Value = "080-01-9053891"

# Numero Propuesta (S2) estaba vacio y debe quedar con el nuevo valor como texto.
# Forzamos formato de texto para que no se interprete como numero, y luego
# restauramos el estilo "Normal" para que la celda no quede con un formato
# numerico aplicado.
$ws.Range("S2").NumberFormat = "@"
$ws.Range("S2").Value = "4873680"
$ws.Range("S2").Style = "Normal"

# Actualiza la celda activa / seleccion visible de la hoja, y desplaza la
# ventana para que la columna I quede como esquina superior izquierda visible.
$win = $excel.ActiveWindow
$win.ScrollColumn = 9
$win.ScrollRow = 1
$ws.Range("G3").Select()
